$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") date values from 45208 (2023-10-09) to
# 45212 (2023-10-13) for rows 2 through 11, keeping existing formatting.
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}
